# Carbohidrates.xlsx -- "Add files via upload" re-save edit
#
# The commit re-uploads the workbook after a short editing session in Excel:
#   1. The header in B1 was renamed from "Value" to "Value (g)" to clarify units.
#   2. Columns A and B were resized to best-fit their (now wider) content.
#   3. The view was scrolled down (row 16 at the top) with B2 as the active cell,
#      instead of the original "whole column A selected" view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clarify the units of the "Value" column ---
$ws.Range("B1").Value = "Value (g)"

# --- 2. Best-fit column widths for A (dates) and B (values) ---
$ws.Columns.Item(1).ColumnWidth = 16.83333
$ws.Columns.Item(2).ColumnWidth = 10.83333

# --- 3. Scroll the window and move the selection to B2 ---
$excel.ActiveWindow.ScrollRow = 16
$null = $ws.Range("B2").Select()
